$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "BOAddon" is being split into two separate order-target columns:
# "BOAddonPersonal" (reusing the old BOAddon column, E) and a new
# "BOAddonRetail" column inserted right after it (shifting BORetail..CCShow
# one column to the right, from F:L to G:M).
$ws.Columns("F:F").Insert()

# Rename the header for the old BOAddon column and label the newly
# inserted column.
$ws.Range("E1").Value = "BOAddonPersonal"
$ws.Range("F1").Value = "BOAddonRetail"

# Populate the new BOAddonRetail column - warehouse shipping is not yet
# enabled for any ordering window.
$ws.Range("F2:F5").Value = "no"

# Match column widths to the new BOAddonPersonal / BOAddonRetail columns
# (values tuned so the engine's rendered ColumnWidth lands on ~18.57 / ~19.14
# characters, matching the widths Excel auto-fit to for these headers).
$ws.Columns("E:E").ColumnWidth = 17.6
$ws.Columns("F:F").ColumnWidth = 18.3

# Reflect the selection left after making these edits.
$ws.Range("F2:F5").Select()
